$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix check for current holdings: rows 3-5 (ZDBDID 2278487-1, 2278491-3, 2606456-X)
# previously showed 0 / blank for number-of-libraries / libraries list;
# they actually have 1 holding library (FL-Bibliothek 547).
# Force text storage (matches source data type) via text number format.
$ws.Range("X3:Y5").NumberFormat = "@"

$ws.Range("X3").Value = "1"
$ws.Range("Y3").Value = "547"

$ws.Range("X4").Value = "1"
$ws.Range("Y4").Value = "547"

$ws.Range("X5").Value = "1"
$ws.Range("Y5").Value = "547"

# Reset formatting so the cells keep the workbook's default (unstyled) look,
# matching the rest of the sheet.
$ws.Range("X3:Y5").ClearFormats()
